# Insert a new weekly price record as row 504, pushing the existing
# rows 504:537 down to 505:538 (dimension grows from A1:R537 to A1:R538).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 504 (shifts rows 504-537 -> 505-538).
$ws.Rows.Item(504).Insert()

# Populate the new row 504 with the added record.
$ws.Range("A504").Value = 6
$ws.Range("B504").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C504").Value = "Metropolitana"
$ws.Range("D504").Value = 44826
$ws.Range("E504").Value = 13
$ws.Range("F504").Value = 100112039
$ws.Range("G504").Value = "Ciboulette"
$ws.Range("H504").Value = "Sin especificar"
$ws.Range("I504").Value = "Primera"
$ws.Range("J504").Value = 830
$ws.Range("K504").Value = 900
$ws.Range("L504").Value = 1000
$ws.Range("M504").Value = 958
$ws.Range("N504").Value = "`$/docena de atados"
$ws.Range("O504").Value = "Región Metropolitana"
$ws.Range("P504").Value = 319
$ws.Range("Q504").Value = 3
$ws.Range("R504").Value = "Hortaliza"
